# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") on Sheet1 held a previously-computed "Strike#" value.
# This regenerates that column from the real strikeout ("K") counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> new K value (column G), keyed by the sheet row number
$kValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 1
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 3
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 3
    19 = 2
    20 = 2
    21 = 1
    22 = 0
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
